$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "owner" -> "editor" rename (owner_email_list -> editor_email_list), to
# match datastore v4.0.0. The header cell on every template sheet that held
# "owner_email_list" is updated to "editor_email_list". Updating every
# reference causes the now-unused shared string to drop out of the workbook's
# shared string table (matching the OOXML diff exactly).
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Project")
$ws.Range("K1").Value = "editor_email_list"
[void]$ws.Range("K1").Select()

$ws = $wb.Worksheets.Item("GenericDataset")
$ws.Range("K1").Value = "editor_email_list"
[void]$ws.Range("K1").Select()

$ws = $wb.Worksheets.Item("FieldNotes")
$ws.Range("K1").Value = "editor_email_list"
[void]$ws.Range("K1").Select()

$ws = $wb.Worksheets.Item("GenericDocument")
$ws.Range("K1").Value = "editor_email_list"
[void]$ws.Range("K1").Select()

$ws = $wb.Worksheets.Item("AudioRecording")
$ws.Range("L1").Value = "editor_email_list"
[void]$ws.Range("L1").Select()

$ws = $wb.Worksheets.Item("Script")
$ws.Range("L1").Value = "editor_email_list"
[void]$ws.Range("L1").Select()

$ws = $wb.Worksheets.Item("WebSite")
$ws.Range("K1").Value = "editor_email_list"
[void]$ws.Range("K1").Select()
